# Apply the "more test runs with simpler rule" update to the Log workbook:
#  - append rows 197-233 of new test-run log data to Tabelle1
#  - widen columns C and H slightly to fit the new content
#  - extend the "log" defined name to cover the new rows
#  - leave the selection on the first newly-added row's area (A208)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")


# --- New log rows (row, timestamp, user, id, lra, depth, rule, steps, model) ---
$rowsData = @(
  ,@(197, 42506.756712962961, "zuendorf", 286817, -22.855555555555501, 100, "DEPTHIGNORE", 301, "input_models/TTC_InputRDG_F.xmi")
  ,@(198, 42510.721400462964, "zuendorf", 8989, 4, 100, "DEFAULT", 5, "input_models/TTC_InputRDG_A.xmi")
  ,@(199, 42510.72146990741, "zuendorf", 6368, 4, 100, "DEPTH", 5, "input_models/TTC_InputRDG_A.xmi")
  ,@(200, 42510.72148148148, "zuendorf", 731, 4, 100, "IGNORE", 5, "input_models/TTC_InputRDG_A.xmi")
  ,@(201, 42510.721574074072, "zuendorf", 8076, 4, 100, "DEPTHIGNORE", 5, "input_models/TTC_InputRDG_A.xmi")
  ,@(202, 42510.721608796295, "zuendorf", 3377, 2.5833333333333299, 100, "DEFAULT", 11, "input_models/TTC_InputRDG_B.xmi")
  ,@(203, 42510.721678240741, "zuendorf", 5415, 1.49999999999999, 100, "DEPTH", 9, "input_models/TTC_InputRDG_B.xmi")
  ,@(204, 42510.721724537034, "zuendorf", 4329, 2.5833333333333299, 100, "IGNORE", 11, "input_models/TTC_InputRDG_B.xmi")
  ,@(205, 42510.72184027778, "zuendorf", 10208, 2.2166666666666601, 100, "DEPTHIGNORE", 9, "input_models/TTC_InputRDG_B.xmi")
  ,@(206, 42510.721886574072, "zuendorf", 4357, -7, 100, "DEFAULT", 31, "input_models/TTC_InputRDG_C.xmi")
  ,@(207, 42510.721944444442, "zuendorf", 4373, 3, 100, "DEPTH", 24, "input_models/TTC_InputRDG_C.xmi")
  ,@(208, 42510.722291666665, "zuendorf", 30384, 5.0357142857142803, 100, "IGNORE", 18, "input_models/TTC_InputRDG_C.xmi")
  ,@(209, 42510.722546296296, "zuendorf", 21913, 5.0357142857142803, 100, "DEPTHIGNORE", 18, "input_models/TTC_InputRDG_C.xmi")
  ,@(210, 42510.722708333335, "zuendorf", 14131, -36, 100, "DEFAULT", 78, "input_models/TTC_InputRDG_D.xmi")
  ,@(211, 42510.72283564815, "zuendorf", 11043, -16.630952380952301, 100, "DEPTH", 65, "input_models/TTC_InputRDG_D.xmi")
  ,@(212, 42510.723032407404, "zuendorf", 16308, -34, 100, "IGNORE", 77, "input_models/TTC_InputRDG_D.xmi")
  ,@(213, 42510.723483796297, "zuendorf", 39263, 9.0309433364988898, 100, "DEPTHIGNORE", 33, "input_models/TTC_InputRDG_D.xmi")
  ,@(214, 42510.723946759259, "zuendorf", 39619, -78, 100, "DEFAULT", 159, "input_models/TTC_InputRDG_E.xmi")
  ,@(215, 42510.724259259259, "zuendorf", 27647, -57.6666666666666, 100, "DEPTH", 147, "input_models/TTC_InputRDG_E.xmi")
  ,@(216, 42511.891956018517, "zuendorf", 100757719, 1.5, 50, "DEPTH", 1, "input_models/TTC_InputRDG_Small1.xmi")
  ,@(217, 42513.554537037038, "zuendorf", 4583, 4, 50, "DEPTH", 5, "input_models/TTC_InputRDG_A.xmi")
  ,@(218, 42513.554745370369, "zuendorf", 5308, 4, 50, "DEPTH", 5, "input_models/TTC_InputRDG_A.xmi")
  ,@(219, 42513.575844907406, "zuendorf", 1731090, 4, 50, "DEPTH", 5, "input_models/TTC_InputRDG_A.xmi")
  ,@(220, 42513.576863425929, "zuendorf", 3910, 4, 50, "DEPTH", 5, "input_models/TTC_InputRDG_A.xmi")
  ,@(221, 42513.61822916667, "zuendorf", 3552, 3, 50, "DEPTH", 4, "input_models/TTC_InputRDG_A.xmi")
  ,@(222, 42513.622824074075, "zuendorf", 384751, 3, 50, "DEPTH", 4, "input_models/TTC_InputRDG_A.xmi")
  ,@(223, 42513.62363425926, "zuendorf", 2859, 1.55, 50, "DEPTH", 4, "input_models/TTC_InputRDG_B.xmi")
  ,@(224, 42513.624120370368, "zuendorf", 25715, 1.9083333333333301, 500, "DEPTH", 6, "input_models/TTC_InputRDG_B.xmi")
  ,@(225, 42513.624456018515, "zuendorf", 3355, 3, 100, "DEFAULT", 4, "input_models/TTC_InputRDG_A.xmi")
  ,@(226, 42513.624502314815, "zuendorf", 4375, 3, 100, "DEPTH", 4, "input_models/TTC_InputRDG_A.xmi")
  ,@(227, 42513.624513888892, "zuendorf", 1116, 3, 100, "IGNORE", 4, "input_models/TTC_InputRDG_A.xmi")
  ,@(228, 42513.625740740739, "zuendorf", 97911, 3, 2000, "DEPTH", 8, "input_models/TTC_InputRDG_B.xmi")
  ,@(229, 42513.627962962964, "zuendorf", 7175, 3, 100, "DEFAULT", 4, "input_models/TTC_InputRDG_A.xmi")
  ,@(230, 42513.628009259257, "zuendorf", 4611, 3, 100, "DEPTH", 4, "input_models/TTC_InputRDG_A.xmi")
  ,@(231, 42513.628020833334, "zuendorf", 1133, 3, 100, "IGNORE", 4, "input_models/TTC_InputRDG_A.xmi")
  ,@(232, 42513.631597222222, "zuendorf", 302556, 3, 5000, "DEPTH", 8, "input_models/TTC_InputRDG_B.xmi")
  ,@(233, 42513.641412037039, "zuendorf", 807858, 3, 10000, "DEPTH", 8, "input_models/TTC_InputRDG_B.xmi")
)

foreach ($row in $rowsData) {
  $r   = $row[0]
  $ws.Cells.Item($r, 1).Value = $row[1]
  $ws.Cells.Item($r, 2).Value = $row[2]
  $ws.Cells.Item($r, 3).Value = $row[3]
  $ws.Cells.Item($r, 4).Value = $row[4]
  $ws.Cells.Item($r, 5).Value = $row[5]
  $ws.Cells.Item($r, 6).Value = $row[6]
  $ws.Cells.Item($r, 7).Value = $row[7]
  $ws.Cells.Item($r, 8).Value = $row[8]
}

# Column A on the new rows needs the same date/time display as the existing log
# rows (row 196) -- copy its number format down instead of re-deriving a style.
$ws.Range("A196").Copy()
$ws.Range("A197:A233").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-apply the timestamp values after the format paste (PasteSpecial only moves
# formatting, but re-setting keeps this block self-contained/order independent).
foreach ($row in $rowsData) {
  $ws.Cells.Item($row[0], 1).Value = $row[1]
}

# Columns C and H need to be a little wider to fit the new, longer values.
$ws.Columns.Item(3).ColumnWidth = 9.14
$ws.Columns.Item(8).ColumnWidth = 37.3

# The "log" named range now spans the appended rows too.
$logName = $ws.Names.Item(1)
$logName.RefersTo = "=Tabelle1!`$A`$1:`$H`$233"

# Match the author's final cursor position on the new data.
[void]$ws.Range("A208").Select()
